$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "total" column (F) = quantite * prix_unit for each data row (2..16)
$ws.Range("F2").Formula = "=+E2*D2"
$ws.Range("F3:F16").Formula = "=+E3*D3"

# Grand total row
$ws.Range("F17").Formula = "=SUM(F2:F16)"

# Apply AutoFilter over the data range
$ws.Range("A1:E16").AutoFilter() | Out-Null

# Register the hidden sheet-scoped "_FilterDatabase" name that Excel
# normally writes alongside an AutoFilter
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Feuil1!`$A`$1:`$E`$16")
$filterName.Visible = $false

# Update the active selection to match the recorded state
$ws.Range("E10").Select() | Out-Null
